# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 122 of the "Achicoria" sheet
# (Vega Modelo de Temuco). All the existing records from the old row 122
# down to the old row 137 shift down by one row (to rows 123-138), and the
# new record carries the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 122 downward one row lower, freeing up row 122
# for the newly reported observation.
$ws.Rows.Item(122).Insert()

$ws.Range("A122").Value = 10
$ws.Range("B122").Value = "Vega Modelo de Temuco"
$ws.Range("C122").Value = "La Araucanía"
$ws.Range("D122").Value = 45154
$ws.Range("E122").Value = 9
$ws.Range("F122").Value = 100112010
$ws.Range("G122").Value = "Achicoria"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 50
$ws.Range("K122").Value = 10000
$ws.Range("L122").Value = 10000
$ws.Range("M122").Value = 10000
$ws.Range("N122").Value = "$/caja 18 unidades"
$ws.Range("O122").Value = "Región Metropolitana"
$ws.Range("P122").Value = 556
$ws.Range("Q122").Value = 18
$ws.Range("R122").Value = "Hortaliza"
